# "time_table polished, minor changes"
# Sheet1 holds a small time-table: class_names, t_1/l_1, t_2/l_2, t_3/l_3 ...
# This edit updates several numeric loads/counts in columns C, E, G and
# introduces a 4th teacher/lesson pair (columns H:I) for every class row,
# plus relabels a couple of teacher codes (B/F columns) for classes 3A-5A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 (class 1A) ---
$ws.Range("C2").Value = 10
$ws.Range("E2").Value = 12
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = "d"
$ws.Range("I2").Value = 3

# --- Row 3 (class 2A) ---
$ws.Range("C3").Value = 12
$ws.Range("E3").Value = 8
$ws.Range("G3").Value = 7
$ws.Range("H3").Value = "d"
$ws.Range("I3").Value = 3

# --- Row 4 (class 3A) ---
$ws.Range("B4").Value = "g"
$ws.Range("C4").Value = 12
$ws.Range("E4").Value = 5
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = "h"
$ws.Range("I4").Value = 6

# --- Row 5 (class 4A) ---
$ws.Range("B5").Value = "g"
$ws.Range("C5").Value = 15
$ws.Range("G5").Value = 7
$ws.Range("H5").Value = "h"
$ws.Range("I5").Value = 6

# --- Row 6 (class 5A) ---
$ws.Range("B6").Value = "g"
$ws.Range("C6").Value = 15
$ws.Range("G6").Value = 7
$ws.Range("H6").Value = "h"
$ws.Range("I6").Value = 6

# --- View: zoom level and current selection moved to G5 ---
$excel.ActiveWindow.Zoom = 205
$ws.Range("G5").Select()
